# Access_Data_Export/clients.xlsx
#
# ImportController fix: the Access export used to hand back the raw lookup
# *labels* ("status", "cash_or_probono", "power_of_attorney_location",
# "documents_location") but the importer actually needs the foreign-key
# *id* columns, so the sheet's header row is renamed to the *_id variants
# and every data row's F/H/I columns are rewritten from the old label-table
# row numbers to the new id-table row numbers (offset by the rows the
# other lookup tables occupy ahead of them: +3 for status_id/cash_or_probono_id,
# +6 for power_of_attorney_location_id, +9 for documents_location_id).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename the four lookup columns to their *_id form ---
# (Table1's ListColumn names follow the header cell text automatically.)
# Order matters for shared-string interning parity with the target file:
# status_id, then cash_or_probono_id, then power_of_attorney_location_id,
# then documents_location_id.
$ws.Range("F1").Value = "status_id"
$ws.Range("E1").Value = "cash_or_probono_id"
$ws.Range("H1").Value = "power_of_attorney_location_id"
$ws.Range("I1").Value = "documents_location_id"

# --- Data rows: renumber the id columns ---
$lastRow = 309

for ($r = 2; $r -le $lastRow; $r++) {
    $fCell = $ws.Range("F" + $r)
    $fVal = $fCell.Value()
    if ($fVal -ne $null) {
        $fCell.Value = $fVal + 3
    }

    $hCell = $ws.Range("H" + $r)
    $hVal = $hCell.Value()
    if ($hVal -ne $null) {
        $hCell.Value = $hVal + 6
    }

    $iCell = $ws.Range("I" + $r)
    $iVal = $iCell.Value()
    if ($iVal -ne $null) {
        $iCell.Value = $iVal + 9
    }
}

# --- Selection tidy-up: collapse the H:I column selection down to I1 ---
$ws.Range("I1").Select()
